$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Add the "_GoBack" bookmark at the very start of the document (the
#    empty first paragraph) - this is what real Word stamps at the last
#    edit position whenever a document is edited & saved.
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$goBackRange = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------
# 2. Collapse the "Diffie-Helman is generally explained ..." paragraph
#    (originally split across 3 runs with proofErr markers around
#    "Alice") down into a single run holding the whole sentence.
# ---------------------------------------------------------------------
$diffieHelman = "Diffie-Helman is generally explained by two sample parties, Alice and Bob, initiating a dialogue."
$rng = $d.Content
$rng.Start = 0
$rng.Find.Execute($diffieHelman, $false, $false, $false, $false, $false, $true, 1, $false, $diffieHelman, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Collapse the "Diffie-Hellman " + "Algorithm" runs (in the
#    Conclusion paragraph) into a single run "Diffie-Hellman Algorithm".
# ---------------------------------------------------------------------
$conclusion = "implement and design of Diffie-Hellman Algorithm."
$rng2 = $d.Content
$rng2.Start = 0
$rng2.Find.Execute($conclusion, $false, $false, $false, $false, $false, $true, 1, $false, $conclusion, 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Remove the VML "watermark" picture from every header (header1 /
#    header2 / header3 each carry one as a floating w:pict Shape). This
#    leaves header1 & header3 as an empty paragraph, and header2 keeps
#    its department-banner InlineShape intact.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$headers = $sec.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $h = $headers.Item($i)
    for ($j = $h.Shapes.Count; $j -ge 1; $j--) {
        $h.Shapes.Item($j).Delete()
    }
}
